$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ G = new_value; H = new_value; I = new_value (only row 2) }
$changes = @{
    2  = @{ G = 0.1109282614759383;  H = 67.26761617893372;   I = -5.829070540232857 }
    3  = @{ G = 0.153288179077933;   H = 29.61647369184415 }
    4  = @{ G = -0.2881273104338656; H = -4.775338745733944 }
    5  = @{ G = -0.3336595871144162; H = 16.37915843823089 }
    6  = @{ G = 0.1862176858016992;  H = -5.544783746654026 }
    7  = @{ G = 0.2933361746628487;  H = 41.44707791015222 }
    8  = @{ G = 0.07969333450842596; H = -21.78954922417756 }
    9  = @{ G = 0.1301092900904414;  H = 2.869407886470347 }
    10 = @{ G = 0.0348079360263512;  H = -43.34443396865068 }
    11 = @{ G = 0.08179399073235775; H = 63.81682998735577 }
    12 = @{ G = 0.1038180429345155;  H = 12.15422990434188 }
    13 = @{ G = 0.1081417088164642;  H = 41.90387891820429 }
    14 = @{ G = 0.1980689796475784;  H = -12.34398151552145 }
    15 = @{ G = 0.2567939686737156;  H = 4.23416752115348 }
    16 = @{ G = 0.1210531513228685;  H = 6.425163732752237 }
    17 = @{ G = 0.130052095579332;   H = -12.95528987475452 }
    18 = @{ G = -0.02167403824249554; H = -142.115296369326 }
    19 = @{ G = 0.07169465151351899; H = 195.9923846896486 }
    20 = @{ G = 0.1240855323053454;  H = 45.87801660801036 }
    21 = @{ G = 0.07588283420200163; H = 15.93313378650229 }
    22 = @{ G = 0.203281806123856;   H = 6.118200395887283 }
    23 = @{ G = 0.2523183989992882;  H = 16.97256988142132 }
    24 = @{ G = -0.02800339258122476; H = -635.8530426540916 }
    25 = @{ G = 0.02664859338819538; H = 214.5815252227584 }
    26 = @{ G = 0.2119075707299745;  H = 3.436997027847126 }
    27 = @{ G = 0.2214821585090581;  H = 14.82628830345011 }
    28 = @{ G = 0.002405485565558058; H = -96.40509113636071 }
    29 = @{ G = 0.1056326141597872;  H = 12.05971202378614 }
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    if ($vals.ContainsKey("I")) {
        $ws.Range("I$row").Value = $vals.I
    }
}
